$wb = $excel.ActiveWorkbook

# Sheet1 Row19
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 703.73334
$ws.Cells.Item(19, 9).Value = 679
$ws.Cells.Item(19, 10).Value = 720.2222
$ws.Cells.Item(19, 11).Value = 679
$ws.Cells.Item(19, 12).Value = 720.2222
$ws.Cells.Item(19, 13).Value = -504
$ws.Cells.Item(19, 14).Value = -1070.2222

# Sheet1 Row74
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(74, 8).Value = 3099.9
$ws.Cells.Item(74, 9).Value = 2999.5
$ws.Cells.Item(74, 10).Value = 3125
$ws.Cells.Item(74, 11).Value = 2999.5
$ws.Cells.Item(74, 12).Value = 3125
$ws.Cells.Item(74, 13).Value = -2063.5
$ws.Cells.Item(74, 14).Value = -4997

# Sheet1 Row77
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(77, 8).Value = 3099.9
$ws.Cells.Item(77, 9).Value = 2999.5
$ws.Cells.Item(77, 10).Value = 3125
$ws.Cells.Item(77, 11).Value = 14997.5
$ws.Cells.Item(77, 12).Value = 15625
$ws.Cells.Item(77, 13).Value = -10317.5
$ws.Cells.Item(77, 14).Value = -24985

# Sheet1 Row111
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(111, 8).Value = 3266.2727
$ws.Cells.Item(111, 9).Value = 3021.5
$ws.Cells.Item(111, 10).Value = 3560
$ws.Cells.Item(111, 11).Value = 9064.5
$ws.Cells.Item(111, 12).Value = 10680
$ws.Cells.Item(111, 13).Value = -5997.5
$ws.Cells.Item(111, 14).Value = -16814

# Sheet1 Row113
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(113, 8).Value = 2960.2727
$ws.Cells.Item(113, 9).Value = 2866.8572
$ws.Cells.Item(113, 10).Value = 3123.75
$ws.Cells.Item(113, 11).Value = 2866.8572
$ws.Cells.Item(113, 12).Value = 3123.75
$ws.Cells.Item(113, 13).Value = 387.1428000000001
$ws.Cells.Item(113, 14).Value = -9631.75

# Sheet1 Row129
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(129, 8).Value = 898.35187
$ws.Cells.Item(129, 9).Value = 332.92307
$ws.Cells.Item(129, 10).Value = 1077.6342
$ws.Cells.Item(129, 11).Value = 998.7692099999999
$ws.Cells.Item(129, 12).Value = 3232.9026
$ws.Cells.Item(129, 13).Value = 4001.23079
$ws.Cells.Item(129, 14).Value = -13232.9026

# Sheet1 Row132
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(132, 8).Value = 5751424.5
$ws.Cells.Item(132, 9).Value = 7939907.5
$ws.Cells.Item(132, 10).Value = 6656.875
$ws.Cells.Item(132, 11).Value = 23819722.5
$ws.Cells.Item(132, 12).Value = 19970.625
$ws.Cells.Item(132, 13).Value = -23817192.5
$ws.Cells.Item(132, 14).Value = -25030.625

# Sheet1 Row135
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(135, 8).Value = 30303664
$ws.Cells.Item(135, 10).Value = 200002800
$ws.Cells.Item(135, 12).Value = 1800025200
$ws.Cells.Item(135, 14).Value = -1800030270

# Sheet1 Row138
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(138, 8).Value = 1079.2529
$ws.Cells.Item(138, 9).Value = 538.2766
$ws.Cells.Item(138, 11).Value = 1614.8298
$ws.Cells.Item(138, 13).Value = 3525.1702

# Sheet2 Row2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 8).Value = 50904
$ws.Cells.Item(2, 9).Value = 0
$ws.Cells.Item(2, 11).Value = 0
$ws.Cells.Item(2, 13).Value = $null

# Sheet2 Row45
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(45, 8).Value = 1513.4546
$ws.Cells.Item(45, 9).Value = 1538.6666
$ws.Cells.Item(45, 10).Value = 1400
$ws.Cells.Item(45, 11).Value = 1538.6666
$ws.Cells.Item(45, 12).Value = 1400
$ws.Cells.Item(45, 13).Value = -1161.6666
$ws.Cells.Item(45, 14).Value = -2154

# Sheet2 Row63
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(63, 8).Value = 28573494
$ws.Cells.Item(63, 9).Value = 2009.7391
$ws.Cells.Item(63, 10).Value = 83335500
$ws.Cells.Item(63, 11).Value = 2009.7391
$ws.Cells.Item(63, 12).Value = 83335500
$ws.Cells.Item(63, 13).Value = -1323.7391
$ws.Cells.Item(63, 14).Value = -83336872

# Sheet2 Row66
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(66, 8).Value = 28573494
$ws.Cells.Item(66, 9).Value = 2009.7391
$ws.Cells.Item(66, 10).Value = 83335500
$ws.Cells.Item(66, 11).Value = 10048.6955
$ws.Cells.Item(66, 12).Value = 416677500
$ws.Cells.Item(66, 13).Value = -6616.6955
$ws.Cells.Item(66, 14).Value = -416684364

# Sheet2 Row88
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(88, 8).Value = 2594.2856
$ws.Cells.Item(88, 9).Value = 1000
$ws.Cells.Item(88, 10).Value = 2860
$ws.Cells.Item(88, 11).Value = 1000
$ws.Cells.Item(88, 12).Value = 2860
$ws.Cells.Item(88, 13).Value = -594
$ws.Cells.Item(88, 14).Value = -3672

# Sheet2 Row91
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(91, 8).Value = 2594.2856
$ws.Cells.Item(91, 9).Value = 1000
$ws.Cells.Item(91, 10).Value = 2860
$ws.Cells.Item(91, 11).Value = 1000
$ws.Cells.Item(91, 12).Value = 2860
$ws.Cells.Item(91, 13).Value = 404
$ws.Cells.Item(91, 14).Value = -5668

# Sheet2 Row110
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(110, 8).Value = 2159.3076
$ws.Cells.Item(110, 9).Value = 1639
$ws.Cells.Item(110, 11).Value = 1639
$ws.Cells.Item(110, 13).Value = 406

# Sheet2 Row116
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(116, 8).Value = 50904
$ws.Cells.Item(116, 9).Value = 0
$ws.Cells.Item(116, 11).Value = 0
$ws.Cells.Item(116, 13).Value = $null

# Sheet3 Row3
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 8).Value = 50904
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 13).Value = $null

# Sheet3 Row86
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 3293.75
$ws.Cells.Item(86, 9).Value = 3439.7058
$ws.Cells.Item(86, 10).Value = 2466.6667
$ws.Cells.Item(86, 11).Value = 3439.7058
$ws.Cells.Item(86, 12).Value = 2466.6667
$ws.Cells.Item(86, 13).Value = -2316.7058
$ws.Cells.Item(86, 14).Value = -4712.6667

# Sheet3 Row89
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(89, 8).Value = 3293.75
$ws.Cells.Item(89, 9).Value = 3439.7058
$ws.Cells.Item(89, 10).Value = 2466.6667
$ws.Cells.Item(89, 11).Value = 17198.529
$ws.Cells.Item(89, 12).Value = 12333.3335
$ws.Cells.Item(89, 13).Value = -11582.529
$ws.Cells.Item(89, 14).Value = -23565.3335

# Sheet3 Row134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 3070.4426
$ws.Cells.Item(134, 9).Value = 863.7
$ws.Cells.Item(134, 11).Value = 2591.1
$ws.Cells.Item(134, 13).Value = -56.10000000000036

# Sheet4 Row31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 2174.8333
$ws.Cells.Item(31, 9).Value = 2411.4119
$ws.Cells.Item(31, 10).Value = 1600.2858
$ws.Cells.Item(31, 11).Value = 2411.4119
$ws.Cells.Item(31, 12).Value = 1600.2858
$ws.Cells.Item(31, 13).Value = -2116.4119
$ws.Cells.Item(31, 14).Value = -2190.2858

# Sheet4 Row34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 2174.8333
$ws.Cells.Item(34, 9).Value = 2411.4119
$ws.Cells.Item(34, 10).Value = 1600.2858
$ws.Cells.Item(34, 11).Value = 2411.4119
$ws.Cells.Item(34, 12).Value = 1600.2858
$ws.Cells.Item(34, 13).Value = -2209.4119
$ws.Cells.Item(34, 14).Value = -2004.2858

# Sheet4 Row58
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(58, 8).Value = 769.7857
$ws.Cells.Item(58, 9).Value = 670.7895
$ws.Cells.Item(58, 11).Value = 670.7895
$ws.Cells.Item(58, 13).Value = -467.7895

# Sheet4 Row107
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(107, 8).Value = 775.44446
$ws.Cells.Item(107, 9).Value = 595.8
$ws.Cells.Item(107, 11).Value = 595.8
$ws.Cells.Item(107, 13).Value = 1324.2

# Sheet4 Row122
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(122, 8).Value = 1208.75
$ws.Cells.Item(122, 10).Value = 1400
$ws.Cells.Item(122, 12).Value = 4200
$ws.Cells.Item(122, 14).Value = -9100

# Sheet4 Row136
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(136, 8).Value = 769.7857
$ws.Cells.Item(136, 9).Value = 670.7895
$ws.Cells.Item(136, 11).Value = 2012.3685
$ws.Cells.Item(136, 13).Value = 537.6315

# Sheet5 Row4
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(4, 8).Value = 991066.5
$ws.Cells.Item(4, 9).Value = 833380.3
$ws.Cells.Item(4, 10).Value = 1201314.6
$ws.Cells.Item(4, 11).Value = 2500140.9
$ws.Cells.Item(4, 12).Value = 3603943.8
$ws.Cells.Item(4, 13).Value = -2500028.9
$ws.Cells.Item(4, 14).Value = -3604167.8

# Sheet5 Row12
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(12, 8).Value = 93.51613
$ws.Cells.Item(12, 9).Value = 134
$ws.Cells.Item(12, 10).Value = 76.954544
$ws.Cells.Item(12, 11).Value = 402
$ws.Cells.Item(12, 12).Value = 230.863632
$ws.Cells.Item(12, 13).Value = -229
$ws.Cells.Item(12, 14).Value = -576.8636320000001

# Sheet5 Row14
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(14, 8).Value = 151.55
$ws.Cells.Item(14, 9).Value = 151.55
$ws.Cells.Item(14, 11).Value = 454.65
$ws.Cells.Item(14, 13).Value = -281.65

# Sheet6 Row102
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(102, 8).Value = 1223.931
$ws.Cells.Item(102, 9).Value = 1071.76
$ws.Cells.Item(102, 11).Value = 1071.76
$ws.Cells.Item(102, 13).Value = 550.24

# Sheet7 Row22
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = 640.4091
$ws.Cells.Item(22, 9).Value = 438.57144
$ws.Cells.Item(22, 10).Value = 993.625
$ws.Cells.Item(22, 11).Value = 438.57144
$ws.Cells.Item(22, 12).Value = 993.625
$ws.Cells.Item(22, 13).Value = -143.57144
$ws.Cells.Item(22, 14).Value = -1583.625

# Sheet7 Row27
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(27, 8).Value = 640.4091
$ws.Cells.Item(27, 9).Value = 438.57144
$ws.Cells.Item(27, 10).Value = 993.625
$ws.Cells.Item(27, 11).Value = 438.57144
$ws.Cells.Item(27, 12).Value = 993.625
$ws.Cells.Item(27, 13).Value = -331.57144
$ws.Cells.Item(27, 14).Value = -1207.625

# Sheet7 Row40
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 2627.0454
$ws.Cells.Item(40, 9).Value = 1789.5
$ws.Cells.Item(40, 11).Value = 1789.5
$ws.Cells.Item(40, 13).Value = -1653.5

# Sheet7 Row61
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(61, 8).Value = 1833
$ws.Cells.Item(61, 9).Value = 1749.5
$ws.Cells.Item(61, 10).Value = 2000
$ws.Cells.Item(61, 11).Value = 1749.5
$ws.Cells.Item(61, 12).Value = 2000
$ws.Cells.Item(61, 13).Value = -1547.5
$ws.Cells.Item(61, 14).Value = -2404

# Sheet7 Row113
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(113, 8).Value = 1833
$ws.Cells.Item(113, 9).Value = 1749.5
$ws.Cells.Item(113, 10).Value = 2000
$ws.Cells.Item(113, 11).Value = 1749.5
$ws.Cells.Item(113, 12).Value = 2000
$ws.Cells.Item(113, 13).Value = 420.5
$ws.Cells.Item(113, 14).Value = -6340

# Sheet8 Row14
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(14, 8).Value = 351504.5
$ws.Cells.Item(14, 9).Value = 351504.5
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 351504.5
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 14).Value = -351336.5
$ws.Cells.Item(14, 13).Value = $null

# Sheet8 Row107
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(107, 8).Value = 497.14285
$ws.Cells.Item(107, 9).Value = 461.81818
$ws.Cells.Item(107, 10).Value = 626.6667
$ws.Cells.Item(107, 11).Value = 1385.45454
$ws.Cells.Item(107, 12).Value = 1880.0001
$ws.Cells.Item(107, 13).Value = 534.54546
$ws.Cells.Item(107, 14).Value = -5720.0001

# Sheet8 Row113
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(113, 8).Value = 306.9091
$ws.Cells.Item(113, 10).Value = 369.92856
$ws.Cells.Item(113, 12).Value = 1109.78568
$ws.Cells.Item(113, 14).Value = -5449.78568
